$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 188.3
$ws.Range("I28").Value = 196.11111
$ws.Range("J28").Value = 118
$ws.Range("K28").Value = 196.11111
$ws.Range("L28").Value = 118
$ws.Range("M28").Value = 288.88889
$ws.Range("N28").Value = -1088
$ws.Range("H88").Value = 875.75
$ws.Range("I88").Value = 875.75
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 875.75
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -469.75
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 875.75
$ws.Range("I91").Value = 875.75
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 875.75
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = 528.25
$ws.Range("N91").ClearContents()
$ws.Range("H94").Value = 9211.666999999999
$ws.Range("I94").Value = 8499.75
$ws.Range("J94").Value = 9781.200000000001
$ws.Range("K94").Value = 8499.75
$ws.Range("L94").Value = 9781.200000000001
$ws.Range("M94").Value = -8048.75
$ws.Range("N94").Value = -10683.2
$ws.Range("H107").Value = 526.5454999999999
$ws.Range("I107").Value = 425.7143
$ws.Range("K107").Value = 425.7143
$ws.Range("M107").Value = 1494.2857
$ws.Range("H111").Value = 1019.8947
$ws.Range("I111").Value = 752.9375
$ws.Range("J111").Value = 2443.6667
$ws.Range("K111").Value = 2258.8125
$ws.Range("L111").Value = 7331.000100000001
$ws.Range("M111").Value = 808.1875
$ws.Range("N111").Value = -13465.0001
$ws.Range("H116").Value = 6446.561
$ws.Range("I116").Value = 9467.333000000001
$ws.Range("J116").Value = 4703.8076
$ws.Range("K116").Value = 9467.333000000001
$ws.Range("L116").Value = 4703.8076
$ws.Range("M116").Value = -6025.333000000001
$ws.Range("N116").Value = -11587.8076
$ws.Range("H137").Value = 4656.161
$ws.Range("I137").Value = 1025.5714
$ws.Range("J137").Value = 5715.0835
$ws.Range("K137").Value = 3076.7142
$ws.Range("L137").Value = 17145.2505
$ws.Range("M137").Value = -526.7142000000003
$ws.Range("N137").Value = -22245.2505

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1138.5555
$ws.Range("I5").Value = 1446.7142
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 1446.7142
$ws.Range("L5").Value = 60
$ws.Range("M5").Value = -1334.7142
$ws.Range("N5").Value = -284
$ws.Range("H32").Value = 11119.85
$ws.Range("I32").Value = 9681.569
$ws.Range("J32").Value = 19955
$ws.Range("K32").Value = 9681.569
$ws.Range("L32").Value = 19955
$ws.Range("M32").Value = -9394.569
$ws.Range("N32").Value = -20529
$ws.Range("H61").Value = 3829.1667
$ws.Range("I61").Value = 3722.7273
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 3722.7273
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -3510.7273
$ws.Range("N61").Value = -5424
$ws.Range("H74").Value = 33823.066
$ws.Range("I74").Value = 51547.7
$ws.Range("J74").Value = 1596.4546
$ws.Range("K74").Value = 51547.7
$ws.Range("L74").Value = 1596.4546
$ws.Range("M74").Value = -50673.7
$ws.Range("N74").Value = -3344.4546
$ws.Range("H77").Value = 33823.066
$ws.Range("I77").Value = 51547.7
$ws.Range("J77").Value = 1596.4546
$ws.Range("K77").Value = 257738.5
$ws.Range("L77").Value = 7982.273
$ws.Range("M77").Value = -253370.5
$ws.Range("N77").Value = -16718.273
$ws.Range("H110").Value = 892.25
$ws.Range("I110").Value = 608.375
$ws.Range("J110").Value = 1460
$ws.Range("K110").Value = 608.375
$ws.Range("L110").Value = 1460
$ws.Range("M110").Value = 1436.625
$ws.Range("N110").Value = -5550
$ws.Range("H136").Value = 3829.1667
$ws.Range("I136").Value = 3722.7273
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 11168.1819
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -8618.1819
$ws.Range("N136").Value = -20100

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1138.5555
$ws.Range("I4").Value = 1446.7142
$ws.Range("J4").Value = 60
$ws.Range("K4").Value = 1446.7142
$ws.Range("L4").Value = 60
$ws.Range("M4").Value = -1331.7142
$ws.Range("N4").Value = -290
$ws.Range("H22").Value = 1144.037
$ws.Range("I22").Value = 2984.1428
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 2984.1428
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -2811.1428
$ws.Range("N22").Value = -846
$ws.Range("H86").Value = 502806.22
$ws.Range("I86").Value = 1983
$ws.Range("J86").Value = 878423.6
$ws.Range("K86").Value = 1983
$ws.Range("L86").Value = 878423.6
$ws.Range("M86").Value = -860
$ws.Range("N86").Value = -880669.6
$ws.Range("H89").Value = 502806.22
$ws.Range("I89").Value = 1983
$ws.Range("J89").Value = 878423.6
$ws.Range("K89").Value = 9915
$ws.Range("L89").Value = 4392118
$ws.Range("M89").Value = -4299
$ws.Range("N89").Value = -4403350

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 7560
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 9200
$ws.Range("K4").Value = 1000
$ws.Range("L4").Value = 9200
$ws.Range("M4").Value = -888
$ws.Range("N4").Value = -9424
$ws.Range("H31").Value = 23022.78
$ws.Range("I31").Value = 26511.05
$ws.Range("J31").Value = 15679.053
$ws.Range("K31").Value = 26511.05
$ws.Range("L31").Value = 15679.053
$ws.Range("M31").Value = -26216.05
$ws.Range("N31").Value = -16269.053
$ws.Range("H34").Value = 23022.78
$ws.Range("I34").Value = 26511.05
$ws.Range("J34").Value = 15679.053
$ws.Range("K34").Value = 26511.05
$ws.Range("L34").Value = 15679.053
$ws.Range("M34").Value = -26309.05
$ws.Range("N34").Value = -16083.053
$ws.Range("H58").Value = 2647.5435
$ws.Range("I58").Value = 902.72
$ws.Range("J58").Value = 4724.7144
$ws.Range("K58").Value = 902.72
$ws.Range("L58").Value = 4724.7144
$ws.Range("M58").Value = -699.72
$ws.Range("N58").Value = -5130.7144
$ws.Range("H132").Value = 25893.56
$ws.Range("I132").Value = 34676
$ws.Range("J132").Value = 1941.4546
$ws.Range("K132").Value = 104028
$ws.Range("L132").Value = 5824.3638
$ws.Range("M132").Value = -101498
$ws.Range("N132").Value = -10884.3638
$ws.Range("H136").Value = 2647.5435
$ws.Range("I136").Value = 902.72
$ws.Range("J136").Value = 4724.7144
$ws.Range("K136").Value = 2708.16
$ws.Range("L136").Value = 14174.1432
$ws.Range("M136").Value = -158.1599999999999
$ws.Range("N136").Value = -19274.1432

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 439574
$ws.Range("I4").Value = 918332.4
$ws.Range("J4").Value = 712.1667
$ws.Range("K4").Value = 2754997.2
$ws.Range("L4").Value = 2136.5001
$ws.Range("M4").Value = -2754885.2
$ws.Range("N4").Value = -2360.5001
$ws.Range("H117").Value = 1446.1428
$ws.Range("I117").Value = 465
$ws.Range("J117").Value = 1609.6666
$ws.Range("K117").Value = 1395
$ws.Range("L117").Value = 4828.9998
$ws.Range("M117").Value = 2047
$ws.Range("N117").Value = -11712.9998
$ws.Range("H121").Value = 46247.184
$ws.Range("I121").Value = 360
$ws.Range("J121").Value = 56444.332
$ws.Range("K121").Value = 1080
$ws.Range("L121").Value = 169332.996
$ws.Range("M121").Value = 230
$ws.Range("N121").Value = -171952.996
$ws.Range("H122").Value = 673.6429000000001
$ws.Range("I122").Value = 520
$ws.Range("K122").Value = 4680
$ws.Range("M122").Value = -2230
$ws.Range("H131").Value = 113187.08
$ws.Range("I131").Value = 435.625
$ws.Range("J131").Value = 137899.72
$ws.Range("K131").Value = 1306.875
$ws.Range("L131").Value = 413699.16
$ws.Range("M131").Value = 3733.125
$ws.Range("N131").Value = -423779.16

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 13382.857
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 13382.857
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 13382.857
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -13606.857
$ws.Range("H52").Value = 60033
$ws.Range("J52").Value = 60033
$ws.Range("L52").Value = 60033
$ws.Range("N52").Value = -60551
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("N53").ClearContents()
$ws.Range("H113").Value = 3760
$ws.Range("I113").Value = 3600
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 3600
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -1430
$ws.Range("N113").Value = -8340

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 14167
$ws.Range("J2").Value = 14167
$ws.Range("L2").Value = 14167
$ws.Range("N2").Value = -14391
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").ClearContents()
$ws.Range("H136").Value = 304801.03
$ws.Range("I136").Value = 418234.34
$ws.Range("J136").Value = 2312.2222
$ws.Range("K136").Value = 1254703.02
$ws.Range("L136").Value = 6936.6666
$ws.Range("M136").Value = -1252153.02
$ws.Range("N136").Value = -12036.6666

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 675.5
$ws.Range("I2").Value = 351
$ws.Range("K2").Value = 351
$ws.Range("M2").Value = -239
$ws.Range("H4").Value = 804
$ws.Range("I4").Value = 600
$ws.Range("J4").Value = 855
$ws.Range("K4").Value = 600
$ws.Range("L4").Value = 855
$ws.Range("M4").Value = -487
$ws.Range("N4").Value = -1081
$ws.Range("H54").Value = 5000
$ws.Range("I54").Value = 5000
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 5000
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -4480
$ws.Range("N54").ClearContents()
$ws.Range("H58").Value = 60000
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 60000
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 60000
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -60616
$ws.Range("H61").Value = 6100
$ws.Range("I61").Value = 4200
$ws.Range("J61").Value = 8000
$ws.Range("K61").Value = 4200
$ws.Range("L61").Value = 8000
$ws.Range("M61").Value = -3908
$ws.Range("N61").Value = -8584
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H136").Value = 2264219
$ws.Range("I136").Value = 2978246
$ws.Range("J136").Value = 836165.4399999999
$ws.Range("K136").Value = 8934738
$ws.Range("L136").Value = 2508496.32
$ws.Range("M136").Value = -8932188
$ws.Range("N136").Value = -2513596.32

Write-Output "Applied all changes"